# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Valor Mora" values for the oldest period (row 16, periodo 2310)
# and the newest period (row 28, periodo 2210) are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 37333
$ws.Range("F28").Value = 40000
